$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "28.908.89"
$ws.Range("E2").Value = "  -1.67%  "

# Row 3
$ws.Range("D3").Value = "1.833.42"
$ws.Range("E3").Value = "  -1.93%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9995"
$ws.Range("E4").Value = "  -0.11%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.58"
$ws.Range("E5").Value = "  +0.38%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6905"
$ws.Range("E6").Value = "  -1.91%  "

# Row 7
$ws.Range("E7").Value = "  -0.08%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07692"
$ws.Range("E8").Value = "  -2.87%  "

# Row 9
$ws.Range("E9").Value = "  -2.69%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.43"
$ws.Range("E10").Value = "  -4.40%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07799"
$ws.Range("E11").Value = "  -0.51%  "

# Row 12
$ws.Range("D12").Value = "1.831.84"
$ws.Range("E12").Value = "  -3.51%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.079"
$ws.Range("E13").Value = "  -1.76%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "90.57"
$ws.Range("E14").Value = "  -3.47%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6811"
$ws.Range("E15").Value = "  -2.88%  "

# Row 16
$ws.Range("E16").Value = "  -1.11%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008294"
$ws.Range("E17").Value = "  -1.26%  "

# Row 18
$ws.Range("D18").Value = "28.906.20"
$ws.Range("E18").Value = "  -1.98%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "243.29"
$ws.Range("E19").Value = "  -3.63%  "

# Row 20
$ws.Range("D20").Value = "2.077.77"
$ws.Range("E20").Value = "  -3.64%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.72"
$ws.Range("E21").Value = "  -2.98%  "

# Row 22
$ws.Range("E22").Value = "  -0.08%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.470"
$ws.Range("E23").Value = "  -2.60%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9999"
$ws.Range("E24").Value = "  -0.07%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "162.59"
$ws.Range("E25").Value = "  +0.57%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1471"
$ws.Range("E26").Value = "  -5.36%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.800"
$ws.Range("E27").Value = "  -2.34%  "

# Row 28
$ws.Range("E28").Value = "  -3.29%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.544"
$ws.Range("E29").Value = "  +2.50%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.208"
$ws.Range("E30").Value = "  -2.52%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.149"
$ws.Range("E31").Value = "  -2.57%  "

# Row 32
$ws.Range("E32").Value = "  -2.50%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05104"
$ws.Range("E33").Value = "  -3.15%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7641"
$ws.Range("E34").Value = "  +1.69%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.844"
$ws.Range("E35").Value = "  -2.87%  "

# Row 36
$ws.Range("E36").Value = "  -3.42%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.687"
$ws.Range("E37").Value = "  -0.80%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01848"
$ws.Range("E38").Value = "  -1.66%  "

# Row 39
$ws.Range("D39").Value = "1.220.78"
$ws.Range("E39").Value = "  -4.19%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.696"
$ws.Range("E40").Value = "  -2.67%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9385"
$ws.Range("E41").Value = "  +5.27%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "108.16"
$ws.Range("E42").Value = "  -1.04%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9995"
$ws.Range("E43").Value = "  -0.08%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.685"
$ws.Range("E44").Value = "  -5.96%  "

# Row 45
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.683"
$ws.Range("E45").Value = "  +0.49%  "

# Row 46
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000123"
$ws.Range("E46").Value = "  -3.17%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5171"
$ws.Range("E47").Value = "  -0.22%  "

# Row 48
$ws.Range("D48").Value = "1.977.36"
$ws.Range("E48").Value = "  -3.35%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "64.13"
$ws.Range("E49").Value = "  -9.65%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.748"

# Row 51
$ws.Range("B51").Value = "Aptos"
$ws.Range("C51").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.901"
$ws.Range("E51").Value = "  -2.24%  "
